# Added 4wk low sales check: update PO_Forecast values for rows 2-9
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 276
$ws.Range("B3").Value = 260
$ws.Range("B4").Value = 213
$ws.Range("B5").Value = 133
$ws.Range("B6").Value = 54
$ws.Range("B7").Value = 38
$ws.Range("B8").Value = 22
$ws.Range("B9").Value = 7
